$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: build one new market sheet by copying "Denmark" (the last existing
# sheet) and placing the copy right after the given anchor sheet.
# ---------------------------------------------------------------------------

function New-MarketSheet($templateName, $afterName, $newName, $jiraRef, $marketName) {
    $template = $wb.Worksheets.Item($templateName)
    $after = $wb.Worksheets.Item($afterName)
    $template.Copy($null, $after)
    $newSheet = $wb.Worksheets.Item($after.Index + 1)
    $newSheet.Name = $newName

    # Jira reference (B4) then Market name (B2) -- this specific order matches
    # the order new shared strings were appended in the authored workbook.
    $newSheet.Range("B4").Value = $jiraRef
    $newSheet.Range("B2").Value = $marketName

    return $newSheet
}

# ---------------------------------------------------------------------------
# Russia: 3 extra rows vs. the Denmark template (MZX252, MZX253, MX4000, ZX4
# are new weight groups; MX1000/MZX254 shift down).
# ---------------------------------------------------------------------------
$russia = New-MarketSheet "Denmark" "Denmark" "Russia" "NGC-2929/T2897" "Russia Market"

$russia.Rows("27:29").Insert()
$russia.Range("A26").Copy()
$russia.Range("A27:A29").PasteSpecial(-4122)

$russia.Range("A25").Value = "MZX252"
$russia.Range("A26").Value = "MZX253"
$russia.Range("A27").Value = "MX1000"
$russia.Range("A28").Value = "MZX254"
$russia.Range("A29").Value = "MX4000"
$russia.Range("A30").Value = "ZX4"

$russia.Rows("3:4").RowHeight = 28.8

# ---------------------------------------------------------------------------
# Finland: 2 extra rows vs. the Denmark template (MZX252, MZX253, MX4000 are
# new weight groups; MX1000/MZX254 shift down, Black Box/Wg/Panels unchanged).
# ---------------------------------------------------------------------------
$finland = New-MarketSheet -templateName "Denmark" -afterName "Russia" -newName "Finland" -jiraRef "NGC-3130/T2940" -marketName "Finland Market"

$finland.Rows("27:28").Insert()
$finland.Range("A26").Copy()
$finland.Range("A27:A28").PasteSpecial(-4122)

$finland.Range("A25").Value = "MZX252"
$finland.Range("A26").Value = "MZX253"
$finland.Range("A27").Value = "MX1000"
$finland.Range("A28").Value = "MZX254"
$finland.Range("A29").Value = "MX4000"

$finland.Rows("3:4").RowHeight = 28.8

# ---------------------------------------------------------------------------
# Hungary: 2 extra rows vs. the Denmark template (MZX252, MZX253, ZX1, ZX4 are
# new weight groups; MX1000 is replaced by ZX1, MZX254 shifts down).
# ---------------------------------------------------------------------------
$hungary = New-MarketSheet -templateName "Denmark" -afterName "Finland" -newName "Hungary" -jiraRef "NGC-3104/T2989" -marketName "Hungary Market"

$hungary.Rows("27:28").Insert()
$hungary.Range("A26").Copy()
$hungary.Range("A27:A28").PasteSpecial(-4122)

$hungary.Range("A25").Value = "MZX252"
$hungary.Range("A26").Value = "MZX253"
$hungary.Range("A27").Value = "ZX1"
$hungary.Range("A28").Value = "MZX254"
$hungary.Range("A29").Value = "ZX4"

$hungary.Rows("3:4").RowHeight = 28.8

# ---------------------------------------------------------------------------
# Hungary becomes the active/selected tab, matching the authored workbook.
# ---------------------------------------------------------------------------
$hungary.Activate()
